$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their textual formatting
# (values like "0.9990" or "242.04" must not be auto-converted to numbers,
# and the percentage strings must remain literal text with surrounding spaces).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.342.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7135'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.04'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3116'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07740'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.81'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08381'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.886.68'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.238'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7151'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.27'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.327.29'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008304'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.961'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.30'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.124.07'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.19'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.898'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1621'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.93'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.023'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.54'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.411'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.299'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.300'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05201'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.175'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.686'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01866'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.713'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.160.19'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.389'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.022.53'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.807'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.401'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4314'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.054'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.12%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.80%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.924'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7743'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.80%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.46'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8925'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.82'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.66%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9989'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.18%  '

Write-Host "Done applying cryptos update"
